$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1698
    $ws.Range("F3").Value = 7844
    $ws.Range("F5").Value = 242
}
